# Seed ONS weekly estimates
#
# 1. Add a new "seed" worksheet after "chart" (so tab order is data, chart, seed).
# 2. Populate it with a header row (week_ended, week_number, estimated_occurrences)
#    and 62 data rows, with column A formatted as yyyy-mm-dd dates.
# 3. Nudge the frozen-pane / selection state on the "data" sheet.

$wb = $excel.ActiveWorkbook

# --- add the "seed" sheet as the last tab -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$seed = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$seed.Name = "seed"

# --- header row -----------------------------------------------------------
$seed.Cells.Item(1,1).Value = "week_ended"
$seed.Cells.Item(1,2).Value = "week_number"
$seed.Cells.Item(1,3).Value = "estimated_occurrences"

# --- data rows: week-ending date (serial), week number, estimated occurrences
$data = @(
    @(43833,1,12431),
    @(43840,2,12139),
    @(43847,3,11746),
    @(43854,4,10914),
    @(43861,5,11094),
    @(43868,6,10710),
    @(43875,7,10877),
    @(43882,8,10795),
    @(43889,9,10647),
    @(43896,10,10984),
    @(43903,11,10834),
    @(43910,12,11401),
    @(43917,13,13787),
    @(43924,14,17897),
    @(43931,15,22038),
    @(43938,16,20922),
    @(43945,17,18694),
    @(43952,18,15825),
    @(43959,19,13712),
    @(43966,20,11948),
    @(43973,21,11354),
    @(43980,22,10216),
    @(43987,23,9971),
    @(43994,24,9453),
    @(44001,25,9204),
    @(44008,26,9661),
    @(44015,27,8740),
    @(44022,28,8662),
    @(44029,29,8874),
    @(44036,30,9050),
    @(44043,31,8852),
    @(44050,32,9008),
    @(44057,33,10142),
    @(44064,34,8836),
    @(44071,35,8655),
    @(44078,36,9033),
    @(44085,37,9143),
    @(44092,38,9474),
    @(44099,39,9584),
    @(44106,40,9957),
    @(44113,41,10308),
    @(44120,42,10384),
    @(44127,43,11031),
    @(44134,44,11388),
    @(44141,45,11754),
    @(44148,46,12339),
    @(44155,47,12326),
    @(44162,48,12447),
    @(44169,49,12675),
    @(44176,50,13045),
    @(44183,51,13138),
    @(44190,52,13532),
    @(44197,53,15012),
    @(44204,1,16507),
    @(44211,2,18696),
    @(44218,3,19473),
    @(44225,4,18160),
    @(44232,5,15920),
    @(44239,6,13983),
    @(44246,7,13281),
    @(44253,8,11745),
    @(44260,9,10480)
)

$row = 2
foreach ($rec in $data) {
    $seed.Cells.Item($row,1).Value = $rec[0]
    $seed.Cells.Item($row,2).Value = $rec[1]
    $seed.Cells.Item($row,3).Value = $rec[2]
    $row = $row + 1
}

# date formatting for column A (rows 2..63)
$seed.Range("A2:A63").NumberFormat = "yyyy\-mm\-dd;@"

# column widths matching the authored sheet
$seed.Columns.Item(1).ColumnWidth = 12.5703125
$seed.Columns.Item(2).ColumnWidth = 14
$seed.Columns.Item(3).ColumnWidth = 22

# --- seed sheet view: freeze header row, scroll to bottom, select B63 -----
$seed.Activate()
$seed.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$seed.Range("B40").Select()
$seed.Range("B63").Select()

# --- "data" sheet view tweaks ----------------------------------------------
$ws1 = $wb.Worksheets.Item("data")
$ws1.Activate()
$ws1.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("A1:BL22").Select()

# restore original active sheet ("chart", the 2nd tab)
$wb.Worksheets.Item("chart").Activate()
